$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 2: CM 145 -> CM 1515, Test 2 -> Test 1515
$ws.Range("A2").Value = "CM 1515"
$ws.Range("B2").Value = "Test 1515"

# Update existing row 3: CM 146 -> CM 1516, Test 3 -> Test 1516
$ws.Range("A3").Value = "CM 1516"
$ws.Range("B3").Value = "Test 1516"

# Add new row 4
$ws.Range("A4").Value = "CM 1517"
$ws.Range("B4").Value = "Test 1517"
$ws.Range("C4").Value = "Zee"
$ws.Range("D4").Value = "Spanish"
